# Lisätty ammuntaetäisyys, jokaisen kierroksen jälkeen syötetään tulokset -KV
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Taul1")

# Row 9 - Katja's section (R:T) new entry
$ws.Range("R8").Copy()
$ws.Range("R9").PasteSpecial(-4122)
$ws.Range("R9").Value = (Get-Date -Year 2024 -Month 2 -Day 11 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("S9").Value = 3

$ws.Range("T8").Copy()
$ws.Range("T9").PasteSpecial(-4122)
$ws.Range("T9").Value = "Ohjelman aloitus"

# Row 10 - Katja's section (R:T) new entry
$ws.Range("R8").Copy()
$ws.Range("R10").PasteSpecial(-4122)
$ws.Range("R10").Value = (Get-Date -Year 2024 -Month 2 -Day 12 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("S10").Value = 2

$ws.Range("T8").Copy()
$ws.Range("T10").PasteSpecial(-4122)
$ws.Range("T10").Value = "Ohjelman muokkausta"

$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("D25").Select()
